$d = $word.ActiveDocument

# 1) Update Tela_ references to shorter form
$d.Content.Find.Execute("Tela_001", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Tela_1", 2)
$d.Content.Find.Execute("Tela_007", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Tela_7", 2)

# 2) Update history table "Pessoa" cells to add "Davi de Jesus Cruz,"
$d.Content.Find.Execute("Idyl Icaro, Wesley Andrade, Victor Lima", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Idyl Icaro, Davi de Jesus Cruz,Wesley Andrade, Victor Lima", 2)
